# Add a new "AUTO_IDATE()" Irrigation Event section to the Field_Overlay sheet,
# and preserve the original (pre-irrigation) layout as a second sheet named
# "Field_Overlay_NoIrrigation".

$wb = $excel.ActiveWorkbook
$orig = $wb.ActiveSheet            # "Field_Overlay" - currently the only sheet

# 1. Duplicate the original sheet, placing the copy right after the original.
#    The copy will become "Field_Overlay_NoIrrigation" (unchanged content,
#    apart from a couple of field-name tweaks below). The original sheet
#    object keeps the name "Field_Overlay" and gets the new Irrigation block.
$orig.Copy($null, $orig)

$newSheet = $wb.Worksheets.Item(1)      # "Field_Overlay" (to receive irrigation rows)
$noIrrig  = $wb.Worksheets.Item(2)      # "Field_Overlay_NoIrrigation"
$noIrrig.Name = "Field_Overlay_NoIrrigation"

# ---------------------------------------------------------------------------
# 2. Tweak the "Field_Overlay_NoIrrigation" copy: rename the FIELD value to
#    FIELD_ORG (its D2 summary concatenation formula recalculates on its own).
# ---------------------------------------------------------------------------
$noIrrig.Range("D8").Value = "FIELD_ORG"

# Restore view state: NoIrrigation copy is not the selected tab, selection D9
$noIrrig.Range("D9").Select()

# ---------------------------------------------------------------------------
# 3. Insert the new "Irrigation Event" block into "Field_Overlay" between the
#    "Organic Matter Application Event" block (ends row 76) and the
#    "Soil Parameters" block (was row 79, now pushed down to row 85).
# ---------------------------------------------------------------------------
$newSheet.Activate()
$newSheet.Rows("77:82").Insert()

# -- Row 79: section title "Irrigation Event" -------------------------------
$newSheet.Range("A73:D73").Copy()
$newSheet.Range("A79").PasteSpecial(-4122) | Out-Null
$newSheet.Range("A79").Value = "!"
$newSheet.Range("B79").Value = "Irrigation Event"
$newSheet.Rows("79:79").RowHeight = 16.8

# -- Row 80: column headers --------------------------------------------------
$newSheet.Range("A74:J74").Copy()
$newSheet.Range("A80").PasteSpecial(-4122) | Out-Null
$newSheet.Range("A80").Value = "!"
$newSheet.Range("B80").Value = "Dome operator"
$newSheet.Range("C80").Value = "Variable to be modified"
$newSheet.Range("D80").Value = "Value or Function"
$newSheet.Range("E80").Value = "Function arguments"

# -- Row 81: descriptive sub-header row --------------------------------------
$newSheet.Range("A69:H69").Copy()
$newSheet.Range("A81").PasteSpecial(-4122) | Out-Null
$newSheet.Range("G69").Copy()
$newSheet.Range("I81").PasteSpecial(-4122) | Out-Null
$newSheet.Range("M69").Copy()
$newSheet.Range("K81").PasteSpecial(-4122) | Out-Null
$newSheet.Rows("81:81").RowHeight = 72

$newSheet.Range("A81").Value = "!"
$newSheet.Range("C81").Value = "Irrigation date (and other irrigation variables)"
$newSheet.Range("D81").Value = "Fertilizer distribution function"
$newSheet.Range("E81").Value = "Number of irrigation applications"
$newSheet.Range("F81").Value = "Base Temperature(oC)"
$newSheet.Range("G81").Value = "#1 irrigation -Cumulative value of Growing degree-day"
$newSheet.Range("H81").Value = "IRVAL - Irrigation amount, depth of water (mm)"
$newSheet.Range("I81").Value = "#1 irrigation -Cumulative value of Growing degree-day"
$newSheet.Range("J81").Value = "IRVAL - Irrigation amount, depth of water (mm)"
$newSheet.Range("K81").Value = "! Dates depend on planting date"

# -- Row 82: data values ------------------------------------------------------
$newSheet.Range("A76:J76").Copy()
$newSheet.Range("A82").PasteSpecial(-4122) | Out-Null

$newSheet.Range("A82").Value = "&"
$newSheet.Range("B82").Value = "FILL"
$newSheet.Range("C82").Value = "IDATE"
$newSheet.Range("D82").Value = "AUTO_IDATE()"
$newSheet.Range("E82").Value = 2
$newSheet.Range("F82").Value = 5
$newSheet.Range("G82").Value = 400
$newSheet.Range("H82").Value = 50
$newSheet.Range("I82").Value = 160
$newSheet.Range("J82").Value = 40
$newSheet.Range("K82").Value = "! The value is not reasonable, just show the idea of function"

# -- Row 83 / 84: blank spacer rows (borrow formatting only) -----------------
$newSheet.Range("A71:D71").Copy()
$newSheet.Range("A83").PasteSpecial(-4122) | Out-Null

$newSheet.Range("C72:D72").Copy()
$newSheet.Range("C84").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 4. Restore the view state on "Field_Overlay": it is the selected tab, with
#    the default A1 selection.
# ---------------------------------------------------------------------------
$newSheet.Activate()
$newSheet.Range("A1").Select()

$excel.CutCopyMode = $false

Write-Host "Irrigation Event block added; Field_Overlay_NoIrrigation created."
